$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.126.90"
$ws.Range("E2").Value = "  -1.69%  "

$ws.Range("D3").Value = "2.347.58"
$ws.Range("E3").Value = "  -2.85%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.81"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.79%  "

$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0811"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.54%  "

$ws.Range("E12").Value = "  +1.24%  "

$ws.Range("D13").Value = "2.707.48"
$ws.Range("E13").Value = "  -2.97%  "

$ws.Range("E14").Value = "  -3.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.88%  "

$ws.Range("D16").Value = "2.368.84"
$ws.Range("E16").Value = "  -2.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.760"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.84%  "

$ws.Range("D18").Value = "40.070.52"
$ws.Range("E18").Value = "  -1.74%  "

$ws.Range("E19").Value = "  -1.97%  "

$ws.Range("E20").Value = "  -2.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("E24").Value = "  -4.83%  "

$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("E26").Value = "  -2.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.18%  "

$ws.Range("E28").Value = "  -3.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.97%  "

$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("E33").Value = "  -2.46%  "

$ws.Range("E34").Value = "  -2.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0719"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.05%  "

$ws.Range("E36").Value = "  -0.46%  "

$ws.Range("E37").Value = "  -5.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0986"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.19%  "

$ws.Range("D42").Value = "1.968.36"
$ws.Range("E42").Value = "  -0.74%  "

$ws.Range("E43").Value = "  -1.57%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.49%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0266"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.79%  "

$ws.Range("D48").Value = "2.566.57"
$ws.Range("E48").Value = "  -3.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "93.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.55%  "
